$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AI+HPS")

$ws.Range("E19").Value = (Get-Date -Year 2021 -Month 5 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F19").Value = (Get-Date -Year 2021 -Month 5 -Day 24 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G19").Value = "W30"
$ws.Range("H19").Value = "E,V"
$ws.Range("I19").Value = "E,V"

$ws.Activate()
$ws.Range("F24").Select()
